# Add two new columns "I0" (I) and "IF" (J) to the worksheet.
# I0 is a constant 1 for every data row; IF mirrors the existing "IP" (H) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - match styling of existing headers (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data rows 2-37: I = 1 (constant), J = same value as column H (IP)
$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2

    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
